# Edit the Blood & lymphatics Y2 schedule workbook to match the "attendance
# app" re-upload:
#   - Rename the "biochemistry cbl" subject everywhere to
#     "Biochemistry Lab/CBL"
#   - Rename the "pos" subject everywhere to "Parasitology SGD/POS"
#   - Fix a duplicate Session number (Year2/C1 group had two "physiology"
#     rows both tagged Session 1; row 152 should be Session 2)
#   - Widen column C (Subject) so the longer subject names are not clipped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Subject" column is column C, data rows run from row 2 to row 154.
$subjectRange = $ws.Range("C2:C154")

# Whole-cell (not substring) replacements so "parasitology"/"pathology lab"
# etc. are left untouched.
$subjectRange.Replace("biochemistry cbl", "Biochemistry Lab/CBL", -4163, 1, $false)
$subjectRange.Replace("pos", "Parasitology SGD/POS", -4163, 1, $false)

# Data fix: row 152 (Year 2 / C1 / physiology) duplicated Session 1 from
# row 153 - it should be Session 2.
$ws.Cells.Item(152, 4).Value = 2

# Column C needs to be wider to comfortably show the new, longer subject
# names.
$ws.Columns.Item(3).ColumnWidth = 20.5
